$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "redid ptolemy" for 116Cd(d,p): new fit produced a shorter level list
# (53 levels instead of 57) plus revised energies / L / spectroscopic
# factors / errors throughout. Drop the four now-unused trailing rows
# (old rows 55-58) so the sheet shrinks back to A1:E54, then rewrite the
# data block (A=index, B=ENERGY, C=L, D=SPECTROSCOPIC_FACTOR, E=ERROR).
$ws.Rows.Item(55).Resize(4).EntireRow.Delete()

$data = @(
    @(0, 1063.9, 4, 0.1126956974645769, 0.005483515480658945),
    @(1, 1085, 2, 0.2241152910682328, 0.002183827440372549),
    @(2, 1175, 0, 0.02232334056419347, 0.0006144038687392698),
    @(3, 1209.564221739463, 2, 0.01360303665416158, 0.0005318480496363928),
    @(4, 1248, 2, 0.05980508152848407, 0.001109785018054344),
    @(5, 1266.460117430137, 0, 0.0005371386615063451, 0.0001678558317207328),
    @(6, 1292.473979699731, 2, 0.01573908404284468, 0.0006688588172258033),
    @(7, 1304.585912295229, 2, 0.01663672242454893, 0.0007224119571323928),
    @(8, 1317.3, 0, 0.003140089665607161, 0.0002822552508410932),
    @(9, 1339.018418868032, 2, 0.006159145261231409, 0.0004149778722106269),
    @(10, 1358.3, 2, 0.3574992448683941, 0.00291854167373542),
    @(11, 1456.533865675461, 5, 0.1423915023131326, 0.01578142594906256),
    @(12, 1471.462544459845, 3, 0.008832343173097472, 0.0007516887806891464),
    @(13, 1494.44871075181, 5, 0.1467020064439292, 0.009123259107209531),
    @(14, 1515.596556909176, 1, 0.007210205360734137, 0.0001588532277716721),
    @(15, 1538.408632604229, 4, 0.03934791007794066, 0.003934791007794067),
    @(16, 1574, 4, 0.05928654367921919, 0.007469170857224465),
    @(17, 1597, 5, 0.3315318119561868, 0.01384587613738135),
    @(18, 1613.149084845743, 0, 0.00220850793064777, 0.000231367497496433),
    @(19, 1625.344244785309, 2, 0.003718596636069186, 0.0005312280908670265),
    @(20, 1649.346715797556, 0, 0.01762243713181142, 0.0006470907147543973),
    @(21, 1732.134406941666, 0, 0.00198334984842432, 0.0003051307459114338),
    @(22, 1747.835619913894, 4, 0.04741367616777521, 0.004619794088142199),
    @(23, 1772.725731049579, 2, 0.005211610304870546, 0.0003758372816012414),
    @(24, 1785.161835361124, 2, 0.03039905772827027, 0.0009570540129861393),
    @(25, 1795.477475116088, 2, 0.008089490923746204, 0.0005832548314926731),
    @(26, 1805.922029018604, 3, 0.01063909924637024, 0.0008554049645322804),
    @(27, 1818, 0, 0.0509667754600661, 0.001267129776631478),
    @(28, 1840, 2, 0.04978899276743497, 0),
    @(29, 1840, 0, 0.01157104794624787, 0),
    @(30, 1841.298544231364, 2, 0.001819003647809125, 0),
    @(31, 1841.298544231364, 0, 0.003015290156373216, 0),
    @(32, 1851.005074483267, 0, 0.003875156198670626, 0.0003522869271518751),
    @(33, 1865.530752678807, 2, 0.04258180179782398, 0.001395244923432165),
    @(34, 1876, 2, 0.03967881976259829, 0.001429624369420605),
    @(35, 1885.117065551615, 2, 0.009708195856608735, 0.001114946863083916),
    @(36, 1896.361077451411, 2, 0.002365921502383671, 0),
    @(37, 1896.361077451411, 0, 0.002004991328338751, 0),
    @(38, 1911.792281651251, 2, 0.01411352233740401, 0.0007947420345334297),
    @(39, 1924.67520477692, 2, 0.01417646132489715, 0.003971615624095699),
    @(40, 1933.550109340654, 2, 0.01468070724291834, 0.003385720763017084),
    @(41, 1944.234562301843, 2, 0.00903995136124792, 0.001458979491846975),
    @(42, 1959.440632592807, 2, 0.01203939872706242, 0.000976944049451366),
    @(43, 1970.376107067389, 2, 0.01400969684373015, 0.001013097912253214),
    @(44, 1986.766402625133, 2, 0.009054450002693627, 0.0006696217043792715),
    @(45, 1995.02814271466, 1, 0.003077214775464751, 0.0001802713029272675),
    @(46, 2014.856993695082, 1, 0.001706680699809291, 0.0001026535288480353),
    @(47, 2048.833064430493, 2, 0.02759963980852419, 0.000925961941627976),
    @(48, 2073.838281975398, 2, 0.03552571887040871, 0.001081086964779978),
    @(49, 2092.576632909581, 2, 0.02221179704774109, 0.0008740055826112505),
    @(50, 2113.24, 3, 0.06822945443892815, 0.002274315147964273),
    @(51, 2133.721396755173, 3, 0.02991686676381232, 0.001736246731828393),
    @(52, 2148.119842376009, 2, 0.009763655430046462, 0.0006571691154838964)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt 5; $j++) {
        $ws.Cells.Item($r, 1 + $j).Value = $vals[$j]
    }
}
